$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.716.81"
$ws.Range("E2").Value = "  +7.09%  "
$ws.Range("D3").Value = "1.740.81"
$ws.Range("E3").Value = "  +4.00%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'334.60"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.3742"
$ws.Range("E7").Value = "  +2.42%  "
$ws.Range("D8").Value = "'48.43"
$ws.Range("E8").Value = "  +4.18%  "
$ws.Range("D9").Value = "'0.3383"
$ws.Range("E9").Value = "  +4.41%  "
$ws.Range("D10").Value = "'1.189"
$ws.Range("E10").Value = "  +4.09%  "
$ws.Range("D11").Value = "'0.07469"
$ws.Range("E11").Value = "  +5.78%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "'6.398"
$ws.Range("E13").Value = "  +5.48%  "
$ws.Range("D14").Value = "'20.42"
$ws.Range("E14").Value = "  +4.31%  "
$ws.Range("D15").Value = "'7.066"
$ws.Range("E15").Value = "  +6.80%  "
$ws.Range("D16").Value = "1.740.83"
$ws.Range("E16").Value = "  +3.87%  "
$ws.Range("D17").Value = "'0.00001079"
$ws.Range("D18").Value = "'0.06723"
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("D19").Value = "'82.65"
$ws.Range("E19").Value = "  +5.01%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "'16.69"
$ws.Range("E21").Value = "  +5.20%  "
$ws.Range("D22").Value = "'6.225"
$ws.Range("E22").Value = "  +5.16%  "
$ws.Range("D23").Value = "'12.76"
$ws.Range("E23").Value = "  -1.54%  "
$ws.Range("D24").Value = "26.720.94"
$ws.Range("E24").Value = "  +7.13%  "
$ws.Range("D25").Value = "'2.468"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").Value = "'1.468"
$ws.Range("E26").Value = "  +24.25%  "
$ws.Range("D27").Value = "'2.413"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("D28").Value = "'152.43"
$ws.Range("E28").Value = "  +2.81%  "
$ws.Range("D29").Value = "'19.65"
$ws.Range("E29").Value = "  +5.15%  "
$ws.Range("D30").Value = "1.938.41"
$ws.Range("E30").Value = "  +4.02%  "
$ws.Range("D31").Value = "'132.32"
$ws.Range("E31").Value = "  +5.69%  "
$ws.Range("D32").Value = "'4.116"
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("D33").Value = "'6.020"
$ws.Range("E33").Value = "  +4.30%  "
$ws.Range("D34").Value = "'0.08603"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").Value = "'1.689"
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("D36").Value = "'12.89"
$ws.Range("E36").Value = "  +4.97%  "
$ws.Range("D37").Value = "'5.425"
$ws.Range("E37").Value = "  +5.02%  "
$ws.Range("D38").Value = "'0.02352"
$ws.Range("E38").Value = "  +4.57%  "
$ws.Range("D39").Value = "'0.2174"
$ws.Range("E39").Value = "  +4.06%  "
$ws.Range("D40").Value = "'0.06259"
$ws.Range("E40").Value = "  +4.09%  "
$ws.Range("D41").Value = "'8.464"
$ws.Range("E41").Value = "  +2.81%  "
$ws.Range("D42").Value = "'1.222"
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("D43").Value = "'0.6251"
$ws.Range("E43").Value = "  +5.14%  "
$ws.Range("D44").Value = "'14.27"
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "'3.924"
$ws.Range("E46").Value = "  +2.11%  "
$ws.Range("D47").Value = "'0.6075"
$ws.Range("E47").Value = "  +5.89%  "
$ws.Range("D48").Value = "'129.06"
$ws.Range("E48").Value = "  +3.25%  "
$ws.Range("D49").Value = "'2.071"
$ws.Range("E49").Value = "  +5.60%  "
$ws.Range("D50").Value = "'0.07216"
$ws.Range("E50").Value = "  +3.21%  "
$ws.Range("D51").Value = "'77.70"
$ws.Range("E51").Value = "  +4.27%  "
